# ADDED SERVICES IN EXCEL SHEET
# Append the owning micro-service tag to each "User Story" description in
# column B, and add a new entry for the "Forgot Password" service in B17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 - Login User story -> tag with [LOGIN SERVICE]
$ws.Range("B8").Value = "Create Login page for a registered user , take required inputs and validate login credentials from DB  . (Same style as home page)  [LOGIN SERVICE]"

# Row 10 - Flight Search story -> tag with [SEARCH SERVICE]
$ws.Range("B10").Value = "Create a flight search page for user to input, source, destination, date, number of passengers ,type of trip from user and show available flights. (route to home page) [SEARCH SERVICE]"

# Row 7 - Register User story -> tag with [CRUD SERVICE]
$ws.Range("B7").Value = "Create Register page for a new user and take required inputs (validate) . (Same style as home page) [CRUD SERVICE]"

# Row 12 - Payment story -> tag with [PAYMENT SERVICE]
$ws.Range("B12").Value = "Create a payment page and method for a user who has selected seat(s) for a flight and wants to book the seats. [PAYMENT SERVICE]"

# Row 13 - Ticket print story -> tag with [PRINT SERVICE]
$ws.Range("B13").Value = "Create a method to print booked tickets [PRINT SERVICE]"

# Row 17 - new cell describing the Forgot Password service
$ws.Range("B17").Value = "FORGOT PASSWORD SERVICE"

# Match the author's final selection/scroll position
$ws.Range("D12").Select()
